# Insert a new data row at row 16 (pushing the existing rows 16-45 down to 17-46),
# matching the weekly "Fruta / hortaliza" update described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 16; formatting/styles of the
# surrounding rows (e.g. the date-formatted column D) are inherited
# automatically by Excel's row-insert behaviour.
$ws.Rows.Item(16).EntireRow.Insert()

# Populate the new row with the new price-report record.
$ws.Cells.Item(16, 1).Value = 10
$ws.Cells.Item(16, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(16, 3).Value = "La Araucanía"
$ws.Cells.Item(16, 4).Value = 44708
$ws.Cells.Item(16, 5).Value = 9
$ws.Cells.Item(16, 6).Value = "Fruta"
$ws.Cells.Item(16, 7).Value = 100107
$ws.Cells.Item(16, 8).Value = "Otros"
$ws.Cells.Item(16, 9).Value = 100107001
$ws.Cells.Item(16, 10).Value = "Caqui"
$ws.Cells.Item(16, 11).Value = "Mankaki"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 45
$ws.Cells.Item(16, 14).Value = 17000
$ws.Cells.Item(16, 15).Value = 18000
$ws.Cells.Item(16, 16).Value = 17444
$ws.Cells.Item(16, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(16, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(16, 19).Value = 969
$ws.Cells.Item(16, 20).Value = 18
